$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 2.36
$ws.Range("C3").Value = 4.77
$ws.Range("D3").Value = 0.36
$ws.Range("E3").Value = 1.17

$ws.Range("E3").Select()
